# Auto-generated script applying 2024-08-07 crime data update
# across Citywide Totals, By Neighborhood, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4773
$ws.Range("K3").Value = 4913
$ws.Range("E4").Value = 2033
$ws.Range("K4").Value = 1014
$ws.Range("K5").Value = 351
$ws.Range("K6").Value = 5532
$ws.Range("E7").Value = 26038
$ws.Range("K7").Value = 16583

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 311
$ws.Range("K7").Value = 1112

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 130
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 361

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 192
$ws.Range("K3").Value = 264
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 198
$ws.Range("K7").Value = 700

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 183
$ws.Range("K6").Value = 166
$ws.Range("K7").Value = 557

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 378

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 115
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 145
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 487
$ws.Range("K8").Value = 1112
$ws.Range("K9").Value = 72
$ws.Range("K11").Value = 323
$ws.Range("K12").Value = 32
$ws.Range("K18").Value = 112
$ws.Range("K19").Value = 498
$ws.Range("K20").Value = 382
$ws.Range("K28").Value = 7
$ws.Range("K29").Value = 883
$ws.Range("K33").Value = 700
$ws.Range("K34").Value = 88
$ws.Range("K36").Value = 214
$ws.Range("K37").Value = 557
$ws.Range("K42").Value = 616
$ws.Range("K43").Value = 147
$ws.Range("K44").Value = 146
$ws.Range("K47").Value = 110
$ws.Range("K48").Value = 208
$ws.Range("K49").Value = 92
$ws.Range("K51").Value = 209
$ws.Range("K52").Value = 434
$ws.Range("K53").Value = 218
$ws.Range("K54").Value = 326
$ws.Range("E63").Value = 368
$ws.Range("K63").Value = 51
$ws.Range("K64").Value = 104
$ws.Range("K65").Value = 378
$ws.Range("K66").Value = 54
$ws.Range("K70").Value = 28
$ws.Range("K71").Value = 54
$ws.Range("K73").Value = 140
$ws.Range("K75").Value = 56
$ws.Range("K78").Value = 195
$ws.Range("K79").Value = 408
$ws.Range("K83").Value = 361
$ws.Range("K84").Value = 123
$ws.Range("K85").Value = 758
$ws.Range("K86").Value = 111
$ws.Range("K89").Value = 233
$ws.Range("K90").Value = 152
$ws.Range("K91").Value = 177
$ws.Range("K94").Value = 217
$ws.Range("K97").Value = 132
$ws.Range("K99").Value = 280
$ws.Range("E101").Value = 26038
$ws.Range("K101").Value = 16583

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 85
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 326

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 253
$ws.Range("K3").Value = 316
$ws.Range("K4").Value = 44
$ws.Range("K6").Value = 245
$ws.Range("K7").Value = 883

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 149
$ws.Range("K3").Value = 156
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 498

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 33
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 192
$ws.Range("K6").Value = 232
$ws.Range("K7").Value = 616

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 134
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 408

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 123
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 382

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 82
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 169
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 487

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 323

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K2").Value = 13
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 54
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 256
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 758

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 118
$ws.Range("K7").Value = 434

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 7
